$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "happy sushi day" expense row (row 9):
# B9 = item name, D9 = amount spent, E9 = receipt marker
$ws.Range("B9").Value = "무모한초밥"
$ws.Range("D9").Value = 31800
$ws.Range("E9").Value = "V"

# Move the active selection to D9 to match the saved view state
$ws.Range("D9").Select()
